$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The source workbook stores every cell as text (inline strings), including
# numeric-looking "Price" values in column D. Plain `Range.Value = "236.69"`
# would be auto-parsed by Excel into a Number (and would silently normalize
# values like "0.1400" -> 0.14, losing the trailing zero), so each target
# D-column cell is pre-formatted as Text ("@") before the literal is written,
# keeping it a text value exactly like the original file.
$textCells = @("D2", "D3", "D4", "D5", "D6", "D7", "D8", "D9", "D10", "D11", "D13", "D14", "D15", "D16", "D17", "D18", "D19", "D20", "D21", "D22", "D23", "D24", "D25", "D26", "D27", "D40", "D41", "D42", "D43", "D44", "D45", "D46", "D47", "D48", "D49")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Row 2
$ws.Range("D2").Value = "236.69"

# Row 3
$ws.Range("D3").Value = "22.07"

# Row 4
$ws.Range("D4").Value = "5.453"

# Row 5
$ws.Range("D5").Value = "0.05633"

# Row 6
$ws.Range("B6").Value = "KuCoinToken"
$ws.Range("C6").Value = "https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs"
$ws.Range("D6").Value = "6.473"
$ws.Range("E6").Value = "5KuCoinTokenKCS"

# Row 7
$ws.Range("B7").Value = "GateToken"
$ws.Range("C7").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D7").Value = "3.365"
$ws.Range("E7").Value = "6GateTokenGT"

# Row 8
$ws.Range("D8").Value = "1.082"

# Row 9
$ws.Range("D9").Value = "0.7873"

# Row 10
$ws.Range("D10").Value = "0.1400"

# Row 11
$ws.Range("D11").Value = "0.07338"

# Row 13
$ws.Range("D13").Value = "0.02980"

# Row 14
$ws.Range("B14").Value = "BitMartToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D14").Value = "0.09243"
$ws.Range("E14").Value = "13BitMartTokenBMX"

# Row 15
$ws.Range("B15").Value = "BitForexToken"
$ws.Range("C15").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D15").Value = "0.001661"
$ws.Range("E15").Value = "14BitForexTokenBF"

# Row 16
$ws.Range("B16").Value = "MCDex"
$ws.Range("C16").Value = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
$ws.Range("D16").Value = "3.254"
$ws.Range("E16").Value = "15MCDexMCB"

# Row 17
$ws.Range("B17").Value = "CoinExToken"
$ws.Range("C17").Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
$ws.Range("D17").Value = "0.04758"
$ws.Range("E17").Value = "16CoinExTokenCET"

# Row 18
$ws.Range("B18").Value = "One"
$ws.Range("C18").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("D18").Value = "0.0005793"
$ws.Range("E18").Value = "17OneONE"

# Row 19
$ws.Range("B19").Value = "TigerCash"
$ws.Range("C19").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D19").Value = "0.006233"
$ws.Range("E19").Value = "18TigerCashTCH"

# Row 20
$ws.Range("B20").Value = "HotbitToken"
$ws.Range("C20").Value = "https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb"
$ws.Range("D20").Value = "0.005100"
$ws.Range("E20").Value = "19HotbitTokenHTB"

# Row 21
$ws.Range("B21").Value = "BitKan"
$ws.Range("C21").Value = "https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan"
$ws.Range("D21").Value = "0.001051"
$ws.Range("E21").Value = "20BitKanKAN"

# Row 22
$ws.Range("B22").Value = "NitroEx"
$ws.Range("C22").Value = "https://coinranking.com/coin/8oiZw6gwYhC+nitroex-ntx"
$ws.Range("D22").Value = "0.0001501"
$ws.Range("E22").Value = "21NitroExNTX"

# Row 23
$ws.Range("B23").Value = "LEO"
$ws.Range("C23").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D23").Value = "3.860"
$ws.Range("E23").Value = "22LEOLEOBestin24h"

# Row 24
$ws.Range("B24").Value = "BTSEToken"
$ws.Range("C24").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("D24").Value = "2.151"
$ws.Range("E24").Value = "23BTSETokenBTSE"

# Row 25
$ws.Range("B25").Value = "BitpandaEcosystemToken"
$ws.Range("C25").Value = "https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best"
$ws.Range("D25").Value = "0.3291"
$ws.Range("E25").Value = "24BitpandaEcosystemTokenBEST"

# Row 26
$ws.Range("B26").Value = "ProBitToken"
$ws.Range("C26").Value = "https://coinranking.com/coin/lQP4d6T2+probittoken-prob"
$ws.Range("D26").Value = "0.1054"
$ws.Range("E26").Value = "25ProBitTokenPROB"

# Row 27
$ws.Range("D27").Value = "0.0004992"

# Row 40
$ws.Range("D40").Value = "0.04111"

# Row 41
$ws.Range("D41").Value = "0.006969"

# Row 42
$ws.Range("D42").Value = "0.003502"
$ws.Range("E42").Value = "41CEJICEJI"

# Row 43
$ws.Range("D43").Value = "0.1036"

# Row 44
$ws.Range("D44").Value = "0.009911"

# Row 45
$ws.Range("D45").Value = "0.00005436"

# Row 46
$ws.Range("D46").Value = "0.00000000750"

# Row 47
$ws.Range("D47").Value = "0.6756"

# Row 48
$ws.Range("D48").Value = "0.03821"

# Row 49
$ws.Range("D49").Value = "0.00002101"
